$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.480.59'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.87%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.866.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.65%  '
$ws.Range("E4").Value = '  -0.52%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.13%  '
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4678'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3734'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07373'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8896'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07947'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.50%  '
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.815.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.428'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.605'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008924'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.513.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.166'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.111.49'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.891'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.089'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.183'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08919'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.49%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7593'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.030'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.173'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.499'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.630'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +11.16%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01976'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.88%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.083'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05284'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.209'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5208'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1651'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.70%  '
$ws.Range("E44").Value = '  +3.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4881'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.50%  '
$ws.Range("E46").Value = '  +2.57%  '
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.657'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06267'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.37%  '
